# Auto-generated edit script: update crypto price/volume columns
# per commit 'Updated symbol list on Tue Jan 17 22:28:15 UTC 2023 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.91"
$ws.Range("E2").Value = "'0.55%"
$ws.Range("D3").Value = "'32.12"
$ws.Range("E3").Value = "'1.69%"
$ws.Range("D4").Value = "'4.989"
$ws.Range("E4").Value = "'-2.86%"
$ws.Range("D5").Value = "'0.07906"
$ws.Range("E5").Value = "'-2.65%"
$ws.Range("D6").Value = "'2.119"
$ws.Range("E6").Value = "'-15.75%"
$ws.Range("D7").Value = "'7.810"
$ws.Range("E7").Value = "'0.18%"
$ws.Range("D8").Value = "'3.798"
$ws.Range("E8").Value = "'-1.99%"
$ws.Range("E9").Value = "'-0.53%"
$ws.Range("D10").Value = "'0.1752"
$ws.Range("E10").Value = "'-0.18%"
$ws.Range("D11").Value = "'0.07962"
$ws.Range("E11").Value = "'8.22%"
$ws.Range("D12").Value = "'0.08654"
$ws.Range("E12").Value = "'-1.34%"
$ws.Range("D13").Value = "'0.03133"
$ws.Range("E13").Value = "'3.49%"
$ws.Range("D14").Value = "'0.1001"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("D15").Value = "'0.001514"
$ws.Range("E15").Value = "'-0.38%"
$ws.Range("D16").Value = "'0.005741"
$ws.Range("E16").Value = "'-4.63%"
$ws.Range("E17").Value = "'2,099.54%"
$ws.Range("D18").Value = "'3.462"
$ws.Range("E18").Value = "'-3.02%"
$ws.Range("D20").Value = "'0.3288"
$ws.Range("E20").Value = "'0.49%"
$ws.Range("D21").Value = "'0.1309"
$ws.Range("E21").Value = "'-2.28%"
$ws.Range("D22").Value = "'4.269"
$ws.Range("E22").Value = "'2.61%"
$ws.Range("D23").Value = "'0.1791"
$ws.Range("E23").Value = "'6.56%"
$ws.Range("D24").Value = "'0.04603"
$ws.Range("E24").Value = "'-0.44%"
$ws.Range("D25").Value = "'0.001238"
$ws.Range("E25").Value = "'-0.21%"
$ws.Range("D26").Value = "'0.004464"
$ws.Range("E26").Value = "'-1.58%"
$ws.Range("E27").Value = "'4.24%"
$ws.Range("D39").Value = "'0.01719"
$ws.Range("E39").Value = "'-2.60%"
$ws.Range("D40").Value = "'0.04779"
$ws.Range("E40").Value = "'3.95%"
$ws.Range("D41").Value = "'0.007433"
$ws.Range("E41").Value = "'7.41%"
$ws.Range("D42").Value = "'0.1357"
$ws.Range("E42").Value = "'-1.24%"
$ws.Range("D43").Value = "'0.002391"
$ws.Range("E43").Value = "'9.21%"
$ws.Range("E44").Value = "'-0.93%"
$ws.Range("D45").Value = "'0.00006002"
$ws.Range("E45").Value = "'-2.45%"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("E47").Value = "'-59.63%"
$ws.Range("E48").Value = "'9.49%"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("E50").Value = "'0.07%"
